$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (46061 -> 46062) for every data row (rows 2 through 95).
$ws.Range("C2:C95").Value = 46062
